# B6-PowerPoint.pptx edit:
#  1) Re-style the three tables (slides 14, 15, 16) from table style
#     {1F1F045A-E792-4012-8FB0-96828C09B910} to
#     {668D83F8-3AC4-49F8-9B47-0450C6D074A4}.
#  2) Swap the deck's theme palette from "Integral / Red Violet" to the
#     stock "Office" palette (the two theme parts exchange their
#     colour schemes).

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$newStyle = "{668D83F8-3AC4-49F8-9B47-0450C6D074A4}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shpIdx = 1; $shpIdx -le $slide.Shapes.Count; $shpIdx++) {
        $shp = $slide.Shapes.Item($shpIdx)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyle)
        }
    }
}

# --- 2) Theme colours --------------------------------------------------
function ToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office theme colour scheme, in ThemeColorScheme slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ToComRGB($officeColors[$i - 1])
}
